# Auto-generated: updates cryptos list (prices / volume%) per commit
# "Updated cryptos list on Mon Dec 11 09:50:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.445.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.97%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.251.75'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -3.58%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.24%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''234.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.68%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -4.01%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''69.84'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -2.29%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.10%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.559'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -3.93%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  +1.09%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''58.50'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +1.13%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''36.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +12.76%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -1.12%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  -4.30%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''2.583.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.71%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''15.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -5.54%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.862'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -3.24%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.248.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -3.60%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''42.259.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -3.09%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.0' + [char]0x2083 + '0979'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.47%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''  -5.15%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''73.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -5.64%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''234.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -6.52%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''2.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +6.37%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E26").Value = '''  -1.77%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -2.81%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -2.04%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -3.29%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''169.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.15%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''20.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -6.38%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.122'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -3.35%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -5.09%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''5.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +2.07%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.0725'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.36%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -6.41%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''3.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -2.28%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''21.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +15.37%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''2.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -3.25%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''6.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -4.91%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +0.68%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''66.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.73%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -8.03%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''8.99'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -2.27%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.103'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.72%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.190'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -2.30%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -0.19%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.0' + [char]0x2083 + '0156'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +27.56%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +11.26%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = '''TrustWalletToken'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''1.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.31%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = '''NEARProtocol'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''2.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -2.63%  '
$ws.Range("E51").Style = "Normal"
